# Article 7: switch the font name used in the example sheet from
# "Comic Sans" to "Open Sans", move the active selection to A2, and widen
# column A to fit the new text (mirrors the authored workbook's column
# resize for this article).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the font-name cell in place (A2 holds the shared string "Comic Sans").
$ws.Range("A2").Value = "Open Sans"

# The saved selection moved to A2 (single cell).
$ws.Range("A2").Select()

# Column A was widened to fit the longer "Open Sans" label (target stored
# OOXML width is ~16.164 characters; 15.33 is the COM-unit input that this
# host's width quantizer maps closest to that stored value).
$ws.Columns.Item(1).ColumnWidth = 15.33
